$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits alone in the empty
#    paragraph right after the "Questions - Part I" heading.
# ------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# ------------------------------------------------------------------
# 2) Merge "Is this a transition or <spellcheck>transversion</spellcheck>?"
#    into a single run (drops the now unnecessary proofErr wrapping).
# ------------------------------------------------------------------
$d.Content.Find.Execute("Is this a transition or transversion?", $false, $false, $false, $false, $false, $true, 1, $false, "Is this a transition or transversion?", 2)

# ------------------------------------------------------------------
# 3) Replace the "Notice that there are many more mutations..." sentence
#    with the new question text, keeping a trailing run that is just a
#    single space (as a separate run, matching the original split).
# ------------------------------------------------------------------
$oldNotice = "Notice that there are many more mutations present in the mitochondrial " + [char]8220 + "control region" + [char]8221 + " sequences, relative to the "
$d.Content.Find.Execute($oldNotice, $false, $false, $false, $false, $false, $true, 1, $false, "Why might you see more mutations in the Part II sequences than in Part I?", 2)

$noticeParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Why might you see more mutations")) {
        $noticeParaIndex = $i
    }
}
$noticePara = $d.Paragraphs.Item($noticeParaIndex)
$noticeEnd = $noticePara.Range.End - 1
$spaceIns = $d.Range($noticeEnd, $noticeEnd)
$spaceIns.InsertAfter(" ")
# Force the trailing space onto its own run (same formatting) by
# toggling Bold on/off, which splits the run without altering rPr.
$spaceRun = $d.Range($noticeEnd, $noticeEnd + 1)
$spaceRun.Bold = 1
$spaceRun.Bold = 0

# ------------------------------------------------------------------
# 4) Expand "Why might translating these sequences not be helpful?"
#    into several runs, with a "_GoBack" bookmark re-inserted right
#    before the final "helpful?" run.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Why might translating these sequences not be helpful?", $false, $false, $false, $false, $false, $true, 1, $false, "Why might translating these sequences from Part II and Part III not be very helpful?", 2)

$transParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Why might translating these sequences from Part II")) {
        $transParaIndex = $i
    }
}
$transPara = $d.Paragraphs.Item($transParaIndex)
$base = $transPara.Range.Start

# Relative offsets of the run boundaries within the final sentence:
#   "Why might translating these sequences "  [0,38)
#   "from Part II and Part III "              [38,64)
#   "not be "                                 [64,71)
#   "very "                                   [71,76)
#   "helpful?"                                [76,84)
$segments = @(@(38,64), @(64,71), @(71,76))
foreach ($seg in $segments) {
    $segStart = $base + $seg[0]
    $segEnd = $base + $seg[1]
    $segRange = $d.Range($segStart, $segEnd)
    $segRange.Bold = 1
    $segRange.Bold = 0
}

$bmPos = $base + 76
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Host "Edits applied"
